$d = $word.ActiveDocument

# Paragraphs (1-based, COM Paragraphs collection, indices as they exist
# BEFORE the later deletion of the two "Usun produkt" / "Zmien ilosc"
# paragraphs in the basket2.php list) that receive strikethrough
# formatting. When $true the paragraph mark itself (pPr/rPr) is also
# struck through (i.e. the whole paragraph range, mark included, was
# selected); when $false only the visible run text gets struck through
# and the paragraph mark is left alone.
$fullMarkStrike = @(85, 86, 87, 88, 89, 92, 96, 98, 99, 100, 101, 104)
$runOnlyStrike  = @(90, 91, 93, 97, 105, 106)

foreach ($i in $fullMarkStrike) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Font.StrikeThrough = 1
}

foreach ($i in $runOnlyStrike) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range
    $trimmed = $d.Range($r.Start, $r.End - 1)
    $trimmed.Font.StrikeThrough = 1
}

# Remove the two bullet points ("Usun produkt" / "Zmien ilosc") from the
# "Strona zamowienia - basket2.php" list; their content is superseded by
# the following two bullets shifting up into their place.
$first = $d.Paragraphs.Item(102)
$last = $d.Paragraphs.Item(103)
$delRange = $d.Range($first.Range.Start, $last.Range.End)
$delRange.Delete()
